# Case with 380 kV done - update recomputed voltage-magnitude (vm_pu)
# results for rows 2-25, columns C,D,F,E,J,K,L,M,N (G stays pinned at 1).
# Note: N25 is intentionally left untouched (matches source diff).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 1.033472160775461
$ws.Range("D2").Value = 1.034309879571508
$ws.Range("E2").Value = 1.041408002600597
$ws.Range("F2").Value = 1.049445165721083
$ws.Range("J2").Value = 1.038596700309459
$ws.Range("K2").Value = 1.037109722035751
$ws.Range("L2").Value = 1.044187602891931
$ws.Range("M2").Value = 1.052202204439152
$ws.Range("N2").Value = 1.016435734895482
$ws.Range("C3").Value = 1.03544866148694
$ws.Range("D3").Value = 1.036092549598916
$ws.Range("E3").Value = 1.043218299993071
$ws.Range("F3").Value = 1.051456291892648
$ws.Range("J3").Value = 1.040210728735579
$ws.Range("K3").Value = 1.038698961400054
$ws.Range("L3").Value = 1.045805899371433
$ws.Range("M3").Value = 1.054022485445754
$ws.Range("N3").Value = 1.017006909062116
$ws.Range("C4").Value = 1.036723925923739
$ws.Range("D4").Value = 1.03724276409273
$ws.Range("E4").Value = 1.044386385570993
$ws.Range("F4").Value = 1.052754473939849
$ws.Range("J4").Value = 1.041251318080538
$ws.Range("K4").Value = 1.039723590024915
$ws.Range("L4").Value = 1.046849325335708
$ws.Range("M4").Value = 1.055196793833904
$ws.Range("N4").Value = 1.017374263863754
$ws.Range("C5").Value = 1.037259193151547
$ws.Range("D5").Value = 1.037725545878715
$ws.Range("E5").Value = 1.044876679306268
$ws.Range("F5").Value = 1.053299497594486
$ws.Range("J5").Value = 1.041687892178773
$ws.Range("K5").Value = 1.040153472794065
$ws.Range("L5").Value = 1.047287109717779
$ws.Range("M5").Value = 1.055689647449765
$ws.Range("N5").Value = 1.017528171583642
$ws.Range("C6").Value = 1.037349017469049
$ws.Range("D6").Value = 1.03780656257878
$ws.Range("E6").Value = 1.044958957242605
$ws.Range("F6").Value = 1.05339096718027
$ws.Range("J6").Value = 1.041761143266605
$ws.Range("K6").Value = 1.040225601453059
$ws.Range("L6").Value = 1.04736056505487
$ws.Range("M6").Value = 1.055772351924051
$ws.Range("N6").Value = 1.017553982574979
$ws.Range("C7").Value = 1.036731081514731
$ws.Range("D7").Value = 1.03724921803903
$ws.Range("E7").Value = 1.044392939896062
$ws.Range("F7").Value = 1.052761759418988
$ws.Range("J7").Value = 1.041257155075226
$ws.Range("K7").Value = 1.039729337537441
$ws.Range("L7").Value = 1.046855178432662
$ws.Range("M7").Value = 1.055203382584191
$ws.Range("N7").Value = 1.017376322453186
$ws.Range("C8").Value = 1.034140901044236
$ws.Range("D8").Value = 1.034913034895387
$ws.Range("E8").Value = 1.042020496078487
$ws.Range("F8").Value = 1.050125501164186
$ws.Range("J8").Value = 1.039142965960823
$ws.Range("K8").Value = 1.037647593373325
$ws.Range("L8").Value = 1.044735295344061
$ws.Range("M8").Value = 1.052818121420896
$ws.Range("N8").Value = 1.016629232106545
$ws.Range("C9").Value = 1.029547567100441
$ws.Range("D9").Value = 1.030770274006141
$ws.Range("E9").Value = 1.037813782196747
$ws.Range("F9").Value = 1.045454937167118
$ws.Range("J9").Value = 1.035387577612422
$ws.Range("K9").Value = 1.033950006735894
$ws.Range("L9").Value = 1.040970454514
$ws.Range("M9").Value = 1.048586977099161
$ws.Range("N9").Value = 1.015295372722479
$ws.Range("C10").Value = 1.026464305950676
$ws.Range("D10").Value = 1.027989623692738
$ws.Range("E10").Value = 1.034990439066692
$ws.Range("F10").Value = 1.042322909519049
$ws.Range("J10").Value = 1.032862672350952
$ws.Range("K10").Value = 1.031464088675047
$ws.Range("L10").Value = 1.038439648565404
$ws.Range("M10").Value = 1.045746068158162
$ws.Range("N10").Value = 1.014394034301038
$ws.Range("C11").Value = 1.025123900402848
$ws.Range("D11").Value = 1.026780828315564
$ws.Range("E11").Value = 1.033763139320059
$ws.Range("F11").Value = 1.040962044892833
$ws.Range("J11").Value = 1.031764036587679
$ws.Range("K11").Value = 1.030382448758078
$ws.Range("L11").Value = 1.037338555146756
$ws.Range("M11").Value = 1.044510850272173
$ws.Range("N11").Value = 1.014000783690352
$ws.Range("C12").Value = 1.024625184340371
$ws.Range("D12").Value = 1.026331089032654
$ws.Range("E12").Value = 1.0333065236029
$ws.Range("F12").Value = 1.040455829196425
$ws.Range("J12").Value = 1.031355129511134
$ws.Range("K12").Value = 1.029979872410066
$ws.Range("L12").Value = 1.036928750016817
$ws.Range("M12").Value = 1.044051245781937
$ws.Range("N12").Value = 1.013854259460923
$ws.Range("C13").Value = 1.024732198639214
$ws.Range("D13").Value = 1.026427593465892
$ws.Range("E13").Value = 1.033404503200311
$ws.Range("F13").Value = 1.04056444765997
$ws.Range("J13").Value = 1.031442879170226
$ws.Range("K13").Value = 1.030066263295587
$ws.Range("L13").Value = 1.037016691637134
$ws.Range("M13").Value = 1.044149868648449
$ws.Range("N13").Value = 1.013885710053951
$ws.Range("C14").Value = 1.025082693470842
$ws.Range("D14").Value = 1.026743667941063
$ws.Range("E14").Value = 1.033725410602492
$ws.Range("F14").Value = 1.040920215997153
$ws.Range("J14").Value = 1.031730253163559
$ws.Range("K14").Value = 1.030349188269273
$ws.Range("L14").Value = 1.03730469718559
$ws.Range("M14").Value = 1.044472875477544
$ws.Range("N14").Value = 1.013988681248463
$ws.Range("C15").Value = 1.025298534092967
$ws.Range("D15").Value = 1.02693831321499
$ws.Range("E15").Value = 1.033923033182023
$ws.Range("F15").Value = 1.04113931905347
$ws.Range("J15").Value = 1.031907203763829
$ws.Range("K15").Value = 1.030523400046155
$ws.Range("L15").Value = 1.037482038878833
$ws.Range("M15").Value = 1.044671785250238
$ws.Range("N15").Value = 1.014052064873853
$ws.Range("C16").Value = 1.026553149464475
$ws.Range("D16").Value = 1.028069745193749
$ws.Range("E16").Value = 1.035071788229111
$ws.Range("F16").Value = 1.042413124769006
$ws.Range("J16").Value = 1.032935470802504
$ws.Range("K16").Value = 1.031535761621134
$ws.Range("L16").Value = 1.03851261219206
$ws.Range("M16").Value = 1.045827936183996
$ws.Range("N16").Value = 1.014420069899895
$ws.Range("C17").Value = 1.027338688704414
$ws.Range("D17").Value = 1.02877817229284
$ws.Range("E17").Value = 1.035791076527278
$ws.Range("F17").Value = 1.043210878276291
$ws.Range("J17").Value = 1.033579030843155
$ws.Range("K17").Value = 1.032169375539015
$ws.Range("L17").Value = 1.039157645108922
$ws.Range("M17").Value = 1.046551779665146
$ws.Range("N17").Value = 1.014650110295421
$ws.Range("C18").Value = 1.02779636782958
$ws.Range("D18").Value = 1.029190928668442
$ws.Range("E18").Value = 1.036210166335453
$ws.Range("F18").Value = 1.043675745109476
$ws.Range("J18").Value = 1.033953895354391
$ws.Range("K18").Value = 1.032538449618426
$ws.Range("L18").Value = 1.039533378106238
$ws.Range("M18").Value = 1.04697349637915
$ws.Range("N18").Value = 1.014784003418759
$ws.Range("C19").Value = 1.027952338491098
$ws.Range("D19").Value = 1.029331591116166
$ws.Range("E19").Value = 1.036352988039169
$ws.Range("F19").Value = 1.043834177304363
$ws.Range("J19").Value = 1.034081628122467
$ws.Range("K19").Value = 1.032664209843082
$ws.Range("L19").Value = 1.039661408591046
$ws.Range("M19").Value = 1.047117208714369
$ws.Range("N19").Value = 1.014829609319113
$ws.Range("C20").Value = 1.027254461006635
$ws.Range("D20").Value = 1.028702212232167
$ws.Range("E20").Value = 1.035713951343409
$ws.Range("F20").Value = 1.04312533354489
$ws.Range("J20").Value = 1.033510036188973
$ws.Range("K20").Value = 1.032101446870572
$ws.Range("L20").Value = 1.039088491457178
$ws.Range("M20").Value = 1.046474168867376
$ws.Range("N20").Value = 1.014625458730945
$ws.Range("C21").Value = 1.024979504615253
$ws.Range("D21").Value = 1.026650612479895
$ws.Range("E21").Value = 1.033630932052991
$ws.Range("F21").Value = 1.040815471462717
$ws.Range("J21").Value = 1.031645651643395
$ws.Range("K21").Value = 1.030265896380903
$ws.Range("L21").Value = 1.037219909277213
$ws.Range("M21").Value = 1.04437778001982
$ws.Range("N21").Value = 1.013958371379994
$ws.Range("C22").Value = 1.023544334081194
$ws.Range("D22").Value = 1.025356403140518
$ws.Range("E22").Value = 1.032316949953445
$ws.Range("F22").Value = 1.039358932156491
$ws.Range("J22").Value = 1.03046865504786
$ws.Range("K22").Value = 1.029107131678038
$ws.Range("L22").Value = 1.036040359360394
$ws.Range("M22").Value = 1.043055116037393
$ws.Range("N22").Value = 1.013536319569406
$ws.Range("C23").Value = 1.024305611183444
$ws.Range("D23").Value = 1.026042902609285
$ws.Range("E23").Value = 1.033013933171654
$ws.Range("F23").Value = 1.040131482123312
$ws.Range("J23").Value = 1.031093064547597
$ws.Range("K23").Value = 1.029721866172883
$ws.Range("L23").Value = 1.036666114214743
$ws.Range("M23").Value = 1.043756728090608
$ws.Range("N23").Value = 1.013760308969049
$ws.Range("C24").Value = 1.02729252146082
$ws.Range("D24").Value = 1.028736536725845
$ws.Range("E24").Value = 1.035748802311413
$ws.Range("F24").Value = 1.043163988914261
$ws.Range("J24").Value = 1.033541213488347
$ws.Range("K24").Value = 1.032132142462993
$ws.Range("L24").Value = 1.03911974057173
$ws.Range("M24").Value = 1.046509239356227
$ws.Range("N24").Value = 1.01463659859396
$ws.Range("C25").Value = 1.030738653274015
$ws.Range("D25").Value = 1.031844499223361
$ws.Range("E25").Value = 1.038904549512156
$ws.Range("F25").Value = 1.04666551374267
$ws.Range("J25").Value = 1.036362103250245
$ws.Range("K25").Value = 1.034909511630666
$ws.Range("L25").Value = 1.041947353936669
$ws.Range("M25").Value = 1.015642309702364

Write-Output "Updated cells"